$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric strings need to be forced to text
# so Excel keeps them as text (matching the original inlineStr "Price" formatting)
# instead of silently converting them into numbers.
$textCells = @("D5", "D6", "D8", "D10", "D13", "D18", "D19", "D21", "D24", "D25", "D26", "D27", "D31", "D36", "D38", "D40", "D41", "D42", "D45", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.421.34"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.805.32"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D5").Value = "227.88"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  +6.87%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "36.24"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "2.064.98"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "11.31"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "1.817.64"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "34.407.31"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "69.99"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "245.58"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "11.50"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +7.96%  "
$ws.Range("D25").Value = "171.69"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "8.10"
$ws.Range("E26").Value = "  +7.70%  "
$ws.Range("D27").Value = "17.49"
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "1.24"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").Value = "1.382.51"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").Value = "0.662"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -10.67%  "
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "82.34"
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.82"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").Value = "0.952"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("D48").Value = "1.966.99"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "103.12"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  -3.38%  "

# Restore default styling on the cells we temporarily formatted as text
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
